# Weekly price update: insert the newest sampling date (2023-07-28, serial 45135)
# as a new Primera/Segunda pair for "Betarraga" at "Vega Monumental Concepción",
# pushing all the existing historical rows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 433 (shifts 433.. down to 435..,
# carrying along the existing formatting, e.g. the date number format on column D).
$ws.Range("A433:R434").EntireRow.Insert()

# ---- New row 433 (Primera) ----
$ws.Cells.Item(433, 1).Value = 11
$ws.Cells.Item(433, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(433, 3).Value = "Bíobío"
$ws.Cells.Item(433, 4).Value = 45135
$ws.Cells.Item(433, 5).Value = 8
$ws.Cells.Item(433, 6).Value = 100114014
$ws.Cells.Item(433, 7).Value = "Betarraga"
$ws.Cells.Item(433, 8).Value = "Sin especificar"
$ws.Cells.Item(433, 9).Value = "Primera"
$ws.Cells.Item(433, 10).Value = 300
$ws.Cells.Item(433, 11).Value = 700
$ws.Cells.Item(433, 12).Value = 700
$ws.Cells.Item(433, 13).Value = 700
$ws.Cells.Item(433, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(433, 15).Value = "Región Metropolitana"
$ws.Cells.Item(433, 16).Value = 140
$ws.Cells.Item(433, 17).Value = 5
$ws.Cells.Item(433, 18).Value = "Hortaliza"

# ---- New row 434 (Segunda) ----
$ws.Cells.Item(434, 1).Value = 11
$ws.Cells.Item(434, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(434, 3).Value = "Bíobío"
$ws.Cells.Item(434, 4).Value = 45135
$ws.Cells.Item(434, 5).Value = 8
$ws.Cells.Item(434, 6).Value = 100114014
$ws.Cells.Item(434, 7).Value = "Betarraga"
$ws.Cells.Item(434, 8).Value = "Sin especificar"
$ws.Cells.Item(434, 9).Value = "Segunda"
$ws.Cells.Item(434, 10).Value = 200
$ws.Cells.Item(434, 11).Value = 500
$ws.Cells.Item(434, 12).Value = 500
$ws.Cells.Item(434, 13).Value = 500
$ws.Cells.Item(434, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(434, 15).Value = "Región Metropolitana"
$ws.Cells.Item(434, 16).Value = 100
$ws.Cells.Item(434, 17).Value = 5
$ws.Cells.Item(434, 18).Value = "Hortaliza"
